# Scheduled-runner market data refresh: update computed price/profit
# columns (H:N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3695.7
$ws.Range("I40").Value = 3695.7
$ws.Range("K40").Value = 3695.7
$ws.Range("M40").Value = -3520.7

$ws.Range("H51").Value = 44216.082
$ws.Range("I51").Value = 2720.3333
$ws.Range("J51").Value = 58048
$ws.Range("K51").Value = 2720.3333
$ws.Range("L51").Value = 58048
$ws.Range("M51").Value = -2236.3333
$ws.Range("N51").Value = -59016

$ws.Range("H61").Value = 1357
$ws.Range("I61").Value = 1446.25
$ws.Range("K61").Value = 4338.75
$ws.Range("M61").Value = -4166.75

$ws.Range("H100").Value = 21000870
$ws.Range("I100").Value = 27946272
$ws.Range("J100").Value = 164663.33
$ws.Range("K100").Value = 27946272
$ws.Range("L100").Value = 164663.33
$ws.Range("M100").Value = -27945731
$ws.Range("N100").Value = -165745.33

$ws.Range("H111").Value = 3481.647
$ws.Range("I111").Value = 3896.5454
$ws.Range("J111").Value = 2721
$ws.Range("K111").Value = 11689.6362
$ws.Range("L111").Value = 8163
$ws.Range("M111").Value = -8622.636200000001
$ws.Range("N111").Value = -14297

$ws.Range("H138").Value = 372241.97
$ws.Range("I138").Value = 1885261.6
$ws.Range("J138").Value = 7030.3447
$ws.Range("K138").Value = 5655784.800000001
$ws.Range("L138").Value = 21091.0341
$ws.Range("M138").Value = -5650644.800000001
$ws.Range("N138").Value = -31371.0341

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3663.1785
$ws.Range("I2").Value = 5495.364
$ws.Range("K2").Value = 5495.364
$ws.Range("M2").Value = -5382.364

$ws.Range("H43").Value = 20675.5
$ws.Range("J43").Value = 20675.5
$ws.Range("L43").Value = 20675.5
$ws.Range("N43").Value = -21301.5

$ws.Range("H97").Value = 9096113
$ws.Range("I97").Value = 5548.95
$ws.Range("K97").Value = 5548.95
$ws.Range("M97").Value = -5052.95

$ws.Range("H106").Value = 46246.668
$ws.Range("J106").Value = 46246.668
$ws.Range("L106").Value = 46246.668
$ws.Range("N106").Value = -48770.668

$ws.Range("H116").Value = 3663.1785
$ws.Range("I116").Value = 5495.364
$ws.Range("K116").Value = 5495.364
$ws.Range("M116").Value = -3201.364

$ws.Range("H122").Value = 670212.2
$ws.Range("I122").Value = 3086.7666
$ws.Range("K122").Value = 9260.299800000001
$ws.Range("M122").Value = -6810.299800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3663.1785
$ws.Range("I3").Value = 5495.364
$ws.Range("K3").Value = 5495.364
$ws.Range("M3").Value = -5381.364

$ws.Range("H20").Value = 2841.697
$ws.Range("I20").Value = 1786.7391
$ws.Range("K20").Value = 1786.7391
$ws.Range("M20").Value = -1539.7391

$ws.Range("H134").Value = 4049.3462
$ws.Range("I134").Value = 3514.9473
$ws.Range("K134").Value = 10544.8419
$ws.Range("M134").Value = -8009.841899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3423
$ws.Range("I31").Value = 2636.5557
$ws.Range("J31").Value = 4838.6
$ws.Range("K31").Value = 2636.5557
$ws.Range("L31").Value = 4838.6
$ws.Range("M31").Value = -2341.5557
$ws.Range("N31").Value = -5428.6

$ws.Range("H34").Value = 3423
$ws.Range("I34").Value = 2636.5557
$ws.Range("J34").Value = 4838.6
$ws.Range("K34").Value = 2636.5557
$ws.Range("L34").Value = 4838.6
$ws.Range("M34").Value = -2434.5557
$ws.Range("N34").Value = -5242.6

$ws.Range("H105").Value = 141884.73
$ws.Range("I105").Value = 176730.92
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 176730.92
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -174983.92
$ws.Range("N105").Value = -5994

$ws.Range("H108").Value = 39109.332
$ws.Range("J108").Value = 51664.5
$ws.Range("L108").Value = 51664.5
$ws.Range("N108").Value = -59344.5

$ws.Range("H141").Value = 590137.8
$ws.Range("J141").Value = 671827.4399999999
$ws.Range("L141").Value = 671827.4399999999
$ws.Range("N141").Value = -682187.4399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 139.8
$ws.Range("I6").Value = 139.8
$ws.Range("K6").Value = 419.4
$ws.Range("M6").Value = -306.4

$ws.Range("H121").Value = 2360.6924
$ws.Range("I121").Value = 1766.4445
$ws.Range("J121").Value = 3697.75
$ws.Range("K121").Value = 5299.333500000001
$ws.Range("L121").Value = 11093.25
$ws.Range("M121").Value = -3989.333500000001
$ws.Range("N121").Value = -13713.25

$ws.Range("H122").Value = 5937.7812
$ws.Range("J122").Value = 7455.7085
$ws.Range("L122").Value = 67101.3765
$ws.Range("N122").Value = -72001.3765

$ws.Range("H129").Value = 1592.8462
$ws.Range("I129").Value = 1334
$ws.Range("J129").Value = 3016.5
$ws.Range("K129").Value = 4002
$ws.Range("L129").Value = 9049.5
$ws.Range("M129").Value = 998
$ws.Range("N129").Value = -19049.5

$ws.Range("H132").Value = 38234.707
$ws.Range("I132").Value = 1525.25
$ws.Range("K132").Value = 13727.25
$ws.Range("M132").Value = -11197.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 5920
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws.Range("H70").Value = 13899.571
$ws.Range("I70").Value = 12199.5
$ws.Range("K70").Value = 12199.5
$ws.Range("M70").Value = -11929.5

$ws.Range("H73").Value = 13899.571
$ws.Range("I73").Value = 12199.5
$ws.Range("K73").Value = 12199.5
$ws.Range("M73").Value = -11263.5

$ws.Range("H97").Value = 23865.445
$ws.Range("I97").Value = 26598.625
$ws.Range("K97").Value = 26598.625
$ws.Range("M97").Value = -26102.625

$ws.Range("H132").Value = 4645.44
$ws.Range("I132").Value = 4749.857
$ws.Range("J132").Value = 4097.25
$ws.Range("K132").Value = 14249.571
$ws.Range("L132").Value = 12291.75
$ws.Range("M132").Value = -11719.571
$ws.Range("N132").Value = -17351.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 39195
$ws.Range("I40").Value = 60789.43
$ws.Range("K40").Value = 60789.43
$ws.Range("M40").Value = -60653.43

$ws.Range("H82").Value = 2710.5
$ws.Range("I82").Value = 2834
$ws.Range("J82").Value = 2657.5715
$ws.Range("K82").Value = 2834
$ws.Range("L82").Value = 2657.5715
$ws.Range("M82").Value = -2473
$ws.Range("N82").Value = -3379.5715

$ws.Range("H85").Value = 2710.5
$ws.Range("I85").Value = 2834
$ws.Range("J85").Value = 2657.5715
$ws.Range("K85").Value = 2834
$ws.Range("L85").Value = 2657.5715
$ws.Range("M85").Value = -1586
$ws.Range("N85").Value = -5153.5715

$ws.Range("H100").Value = 7395.6665
$ws.Range("I100").Value = 3187
$ws.Range("J100").Value = 9500
$ws.Range("K100").Value = 3187
$ws.Range("L100").Value = 9500
$ws.Range("M100").Value = -2646
$ws.Range("N100").Value = -10582

$ws.Range("H132").Value = 484831.25
$ws.Range("I132").Value = 1066619.5
$ws.Range("J132").Value = 5711.5293
$ws.Range("K132").Value = 3199858.5
$ws.Range("L132").Value = 17134.5879
$ws.Range("M132").Value = -3197328.5
$ws.Range("N132").Value = -22194.5879

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 18000
$ws.Range("I52").Value = 18000
$ws.Range("K52").Value = 18000
$ws.Range("M52").Value = -17774

$ws.Range("H100").Value = 27010.785
$ws.Range("I100").Value = 6631.625
$ws.Range("J100").Value = 54183
$ws.Range("K100").Value = 13263.25
$ws.Range("L100").Value = 108366
$ws.Range("M100").Value = -12722.25
$ws.Range("N100").Value = -109448

$ws.Range("H107").Value = 10751.473
$ws.Range("J107").Value = 13154.042
$ws.Range("L107").Value = 39462.126
$ws.Range("N107").Value = -43302.126
